$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text so numeric-looking strings
# (e.g. "0.678", "239.12") are not auto-converted to numbers by Excel,
# matching the original inlineStr/text storage of column D and E.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "43.963.31"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "2.355.51"
$ws.Range("E3").Value = "  -0.20%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "0.678"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").Value = "239.12"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "73.69"
$ws.Range("E7").Value = "  +1.21%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  +9.54%  "

$ws.Range("D10").Value = "0.100"
$ws.Range("E10").Value = "  -0.31%  "

$ws.Range("D11").Value = "57.29"
$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").Value = "32.06"
$ws.Range("E12").Value = "  +9.44%  "

$ws.Range("D13").Value = "7.28"
$ws.Range("E13").Value = "  +8.82%  "

$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "2.708.77"
$ws.Range("E15").Value = "  -0.13%  "

$ws.Range("D16").Value = "16.58"
$ws.Range("E16").Value = "  -1.51%  "

$ws.Range("E17").Value = "  -1.13%  "

$ws.Range("D18").Value = "2.340.72"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").Value = "43.865.16"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  -0.94%  "

$ws.Range("D21").Value = "6.73"
$ws.Range("E21").Value = "  +4.19%  "

$ws.Range("D22").Value = "76.77"
$ws.Range("E22").Value = "  -1.59%  "

$ws.Range("D23").Value = "257.32"
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("D24").Value = "1.93"
$ws.Range("E24").Value = "  +21.49%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("B26").Value = "WEMIXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D26").Value = "3.65"
$ws.Range("E26").Value = "  -2.91%  "

$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.49"
$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("D28").Value = "10.74"
$ws.Range("E28").Value = "  +2.24%  "

$ws.Range("E29").Value = "  +1.38%  "

$ws.Range("D30").Value = "22.63"
$ws.Range("E30").Value = "  +0.78%  "

$ws.Range("D31").Value = "175.33"
$ws.Range("E31").Value = "  +1.48%  "

$ws.Range("D32").Value = "0.128"
$ws.Range("E32").Value = "  -2.98%  "

$ws.Range("D33").Value = "0.135"
$ws.Range("E33").Value = "  +2.48%  "

$ws.Range("D34").Value = "0.0763"
$ws.Range("E34").Value = "  +4.32%  "

$ws.Range("D35").Value = "5.22"
$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "5.45"
$ws.Range("E36").Value = "  +3.41%  "

$ws.Range("D37").Value = "3.75"
$ws.Range("E37").Value = "  -5.02%  "

$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -3.87%  "

$ws.Range("D39").Value = "6.28"
$ws.Range("E39").Value = "  -2.51%  "

$ws.Range("D40").Value = "0.0278"
$ws.Range("E40").Value = "  +3.03%  "

$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  +12.35%  "

$ws.Range("D42").Value = "0.206"
$ws.Range("E42").Value = "  +13.64%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "9.04"
$ws.Range("E43").Value = "  +2.09%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "18.92"
$ws.Range("E44").Value = "  -3.33%  "

$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").Value = "4.78"
$ws.Range("E46").Value = "  +7.09%  "

$ws.Range("B47").Value = "MultiversX"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D47").Value = "58.38"
$ws.Range("E47").Value = "  +10.97%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.51"
$ws.Range("E48").Value = "  +7.07%  "

$ws.Range("D49").Value = "1.24"
$ws.Range("E49").Value = "  -0.54%  "

$ws.Range("D50").Value = "1.17"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").Value = "99.97"
$ws.Range("E51").Value = "  +1.38%  "

# Restore the default (unstyled) cell style now that values are written,
# so no stray style index is left referenced on these cells.
$priceVolRange.Style = "Normal"